$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"); copy the formatting of the existing
# header cell H1 so the new headers reuse the same bold/border/center-top style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for column I (I0), rows 2-49
$iValues = @(8,9,6,8,7,8,2,8,13,9,6,5,8,10,9,7,9,7,6,6,8,6,4,8,7,8,8,6,6,9,6,6,8,7,3,6,9,8,7,7,9,7,5,7,6,4,6,5)

# New data for column J (IF), rows 2-49
$jValues = @(8,9,7,8,7,8,3,8,13,9,6,5,8,10,9,7,9,8,6,6,8,6,5,8,7,8,8,6,6,9,6,7,8,8,4,8,9,8,8,7,9,7,6,7,7,4,6,5)

for ($k = 0; $k -lt $iValues.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}

